# Commit: "Update examples. Set of small tweaks all intended to improve the examples"
# Insert two new worksheets - studyDesignArms and studyDesignEpochs - right after
# the studyDesign sheet (and before mainTimeline), populate them with data, and
# make studyDesignEpochs the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Cells we borrow existing cell formatting from (so the new sheets pick up the
# same styles already used elsewhere in the workbook rather than creating new
# duplicate style records).
$headerFmtSrc = $wb.Worksheets.Item("studyDesignEncounters").Range("B1")   # bold/shaded header style
$bodyFmtSrc   = $wb.Worksheets.Item("studyDesignEncounters").Range("B2")   # plain left/top aligned body style
$blankFmtSrc  = $wb.Worksheets.Item("study").Range("C1")                  # blank left/top aligned style

# ---------------------------------------------------------------------------
# New sheet: studyDesignArms (inserted right after "studyDesign")
# ---------------------------------------------------------------------------
$studyDesign = $wb.Worksheets.Item("studyDesign")
$arms = $wb.Worksheets.Add($null, $studyDesign)
$arms.Name = "studyDesignArms"

$headerFmtSrc.Copy()
$arms.Range("A1:E1").PasteSpecial(-4122)
$bodyFmtSrc.Copy()
$arms.Range("A2:E3").PasteSpecial(-4122)
$blankFmtSrc.Copy()
$arms.Range("E4").PasteSpecial(-4122)

$arms.Cells.Item(1,1).Value = "studyArmName"
$arms.Cells.Item(1,2).Value = "studyArmDescription"
$arms.Cells.Item(1,3).Value = "studyArmType"
$arms.Cells.Item(1,4).Value = "studyArmDataOriginDescription"
$arms.Cells.Item(1,5).Value = "studyArmDataOriginType"

$arms.Cells.Item(2,1).Value = "Active"
$arms.Cells.Item(2,2).Value = "Active Substance"
$arms.Cells.Item(2,3).Value = "Active Comparator Arm"
$arms.Cells.Item(2,4).Value = "Data collected from subjects"
$arms.Cells.Item(2,5).Value = "Data Generated Within Study"

$arms.Cells.Item(3,1).Value = "Placebo"
$arms.Cells.Item(3,2).Value = "Placebo"
$arms.Cells.Item(3,3).Value = "Placebo Comparator Arm"
$arms.Cells.Item(3,4).Value = "Data collected from subjects"
$arms.Cells.Item(3,5).Value = "Data Generated Within Study"

$arms.Columns.Item(1).ColumnWidth = 17.83203125
$arms.Columns.Item(2).ColumnWidth = 24.83203125
$arms.Columns.Item(3).ColumnWidth = 23.5
$arms.Columns.Item(4).ColumnWidth = 32.1640625
$arms.Columns.Item(5).ColumnWidth = 25.6640625

$arms.Range("A1:E1048576").Select()

# ---------------------------------------------------------------------------
# New sheet: studyDesignEpochs (inserted right after "studyDesignArms")
# ---------------------------------------------------------------------------
$epochs = $wb.Worksheets.Add($null, $arms)
$epochs.Name = "studyDesignEpochs"

$headerFmtSrc.Copy()
$epochs.Range("A1:C1").PasteSpecial(-4122)
$bodyFmtSrc.Copy()
$epochs.Range("A2:C5").PasteSpecial(-4122)

$epochs.Cells.Item(1,1).Value = "studyEpochName"
$epochs.Cells.Item(1,2).Value = "studyEpochDescription"
$epochs.Cells.Item(1,3).Value = "studyEpochType"

$epochs.Cells.Item(2,1).Value = "Screening "
$epochs.Cells.Item(2,2).Value = "Screening Epoch"
$epochs.Cells.Item(2,3).Value = "SCREENING"

$epochs.Cells.Item(3,1).Value = "Treatment"
$epochs.Cells.Item(3,2).Value = "Treatment Epoch"
$epochs.Cells.Item(3,3).Value = "TREATMENT"

$epochs.Cells.Item(4,1).Value = "Follow-Up"
$epochs.Cells.Item(4,2).Value = "Follow-up Epoch"
$epochs.Cells.Item(4,3).Value = "FOLLOW-UP"

$epochs.Cells.Item(5,1).Value = "Baseline"
$epochs.Cells.Item(5,2).Value = "Baseline Epoch"
$epochs.Cells.Item(5,3).Value = "BASELINE"

$epochs.Columns.Item(1).ColumnWidth = 18.33203125
$epochs.Columns.Item(2).ColumnWidth = 27.6640625
$epochs.Columns.Item(3).ColumnWidth = 18.33203125

$epochs.Range("A11").Select()

# studyDesignEpochs becomes the active/selected tab.
$epochs.Select()
